$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1 header value tweaks
$ws.Range("B1").Value = 16
$ws.Range("C1").Value = 20
$ws.Range("D1").Value = 16
$ws.Range("E1").Value = 20

# Row 2: delete D2, add B2
$ws.Range("D2").ClearContents()
$ws.Range("B2").Value = 27.73009143525185

# Row 3: delete B3, update C3
$ws.Range("B3").ClearContents()
$ws.Range("C3").Value = 24.014771804472705

# Update the sheet selection to B1:E3
$ws.Range("B1:E3").Select()
